$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the GitHub Actions data refresh commit.
# Price (column D) cells are temporarily set to Text number format before assigning
# numeric-looking strings (e.g. "0.999", "25.46") so Excel stores them as text,
# matching the source inlineStr cells instead of auto-converting to numbers.
# ClearFormats() afterwards restores the original (default/General) cell formatting
# so no stray style attribute is left behind on the cell.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.206.34"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.479.34"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.61%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +2.65%  "

# Row 10
$ws.Range("E10").Value = "  +0.71%  "

# Row 11
$ws.Range("E11").Value = "  -0.45%  "

# Row 12
$ws.Range("E12").Value = "  +0.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.930.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.46"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.127.79"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.20%  "

# Row 16
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.450.76"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.58"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "350.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.82%  "

# Row 21
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.05"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
$ws.Range("E27").Value = "  +0.56%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.66%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "504.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.53%  "

# Row 33
$ws.Range("E33").Value = "  -0.62%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("E35").Value = "  +0.48%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.13"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.91%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.16"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("E39").Value = "  -0.82%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.67%  "

# Row 42
$ws.Range("E42").Value = "  +0.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.99%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.39%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "142.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.32%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.49"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.92%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0259"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.58%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.514"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("E49").Value = "  +0.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.73%  "

# Row 51
$ws.Range("E51").Value = "  +0.63%  "
